$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $text) {
    # Force Excel to store the value as literal text, not auto-convert
    # numeric-looking strings (e.g. "2.04") into a Number.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '42.809.10'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.529.24'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue $ws.Range('D5') '311.35'
$ws.Range('E5').Value = '  +0.12%  '
Set-TextValue $ws.Range('D6') '101.17'
$ws.Range('E6').Value = '  +2.48%  '
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('E8').Value = '  +0.03%  '
Set-TextValue $ws.Range('D10') '35.78'
$ws.Range('E10').Value = '  -0.23%  '
Set-TextValue $ws.Range('D11') '0.0806'
$ws.Range('E11').Value = '  -0.41%  '
Set-TextValue $ws.Range('D12') '7.35'
$ws.Range('E12').Value = '  -1.52%  '
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').Value = '2.916.81'
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D15') '15.34'
$ws.Range('E15').Value = '  -3.16%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.546.30'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('D18').Value = '42.781.17'
$ws.Range('E19').Value = '  -0.77%  '
Set-TextValue $ws.Range('D20') '12.42'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('E21').Value = '  -0.82%  '
Set-TextValue $ws.Range('D22') '69.91'
$ws.Range('E22').Value = '  +0.68%  '
Set-TextValue $ws.Range('D23') '243.85'
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('E24').Value = '  -1.64%  '
Set-TextValue $ws.Range('D25') '2.04'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').Value = '  +0.05%  '
Set-TextValue $ws.Range('D27') '25.53'
$ws.Range('E27').Value = '  -5.80%  '
$ws.Range('E28').Value = '  -2.55%  '
Set-TextValue $ws.Range('D29') '10.19'
$ws.Range('E29').Value = '  -0.03%  '
Set-TextValue $ws.Range('D30') '38.73'
$ws.Range('E30').Value = '  -3.11%  '
Set-TextValue $ws.Range('D31') '161.67'
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('E33').Value = '  +8.56%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E35').Value = '  -1.02%  '
Set-TextValue $ws.Range('D36') '18.36'
$ws.Range('E36').Value = '  -1.77%  '
Set-TextValue $ws.Range('D37') '3.09'
$ws.Range('E37').Value = '  -6.04%  '
$ws.Range('E38').Value = '  -6.80%  '
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('E41').Value = '  +1.15%  '
Set-TextValue $ws.Range('D42') '22.00'
$ws.Range('E42').Value = '  -3.82%  '
$ws.Range('E43').Value = '  +0.27%  '
Set-TextValue $ws.Range('D44') '3.33'
$ws.Range('E44').Value = '  +3.63%  '
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('D46').Value = '1.993.38'
$ws.Range('E46').Value = '  +0.06%  '
Set-TextValue $ws.Range('D47') '9.03'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').Value = '2.770.57'
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('E49').Value = '  -2.06%  '
Set-TextValue $ws.Range('D50') '79.64'
$ws.Range('E50').Value = '  -2.12%  '
Set-TextValue $ws.Range('D51') '72.54'
$ws.Range('E51').Value = '  -1.93%  '
